$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 136, shifting existing rows 136-217 down to 137-218.
$ws.Rows.Item(136).Insert()

# Populate the newly inserted row 136 with the new record's data.
$ws.Cells.Item(136, 1).Value = 9
$ws.Cells.Item(136, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(136, 3).Value = "Metropolitana"
$ws.Cells.Item(136, 4).Value = 44845
$ws.Cells.Item(136, 5).Value = 13
$ws.Cells.Item(136, 6).Value = "Fruta"
$ws.Cells.Item(136, 7).Value = 100101
$ws.Cells.Item(136, 8).Value = "Berries"
$ws.Cells.Item(136, 9).Value = 100101001
$ws.Cells.Item(136, 10).Value = "Arándano (blue)"
$ws.Cells.Item(136, 11).Value = "Sin especificar"
$ws.Cells.Item(136, 12).Value = "Primera"
$ws.Cells.Item(136, 13).Value = 630
$ws.Cells.Item(136, 14).Value = 11000
$ws.Cells.Item(136, 15).Value = 12000
$ws.Cells.Item(136, 16).Value = 11556
$ws.Cells.Item(136, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(136, 18).Value = "Provincia de Linares"
$ws.Cells.Item(136, 19).Value = 5778
$ws.Cells.Item(136, 20).Value = 2

# Keep the date-formatted number format consistent with the other date cells in column D.
$ws.Cells.Item(136, 4).NumberFormat = $ws.Cells.Item(137, 4).NumberFormat
